# Edit script for Vize.docx - applies the changes described by the commit diff.
#
# Summary of changes:
#  1. "(2022)" -> "(2022" / ", ver.2" / ")"                 (split into 3 runs)
#  2. Remove "euroklíč " in 3 bullet points                 (simple text substitution)
#  3. " + Scénáře" -> " +" / " jejich" / " " / "s" / "cénáře"  (split into 5 runs)

function Replace-RangeWithRuns($d, $Start, $End, $Texts) {
    $runsXml = ""
    foreach ($t in $Texts) {
        $escaped = $t.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $preserve = ""
        if ($t -match "^\s" -or $t -match "\s$" -or $t -eq "") {
            $preserve = ' xml:space="preserve"'
        }
        $runsXml = $runsXml + "<w:r><w:t$preserve>$escaped</w:t></w:r>"
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target = $d.Range($Start, $End)
    $target.InsertXML($xml)
}

$d = $word.ActiveDocument

# --- 1. "(2022)" -> "(2022" + ", ver.2" + ")" ------------------------------
$find1 = $d.Content
$find1.Find.Execute("(2022)", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
Replace-RangeWithRuns $d $find1.Start $find1.End @("(2022", ", ver.2", ")")

# --- 2-4. Remove the stray "euroklíč " before "lokac..." ------------------
$d.Content.Find.Execute("Automatickou aktualizaci existujících euroklíč lokací.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Automatickou aktualizaci existujících lokací.", 2) | Out-Null

$d.Content.Find.Execute("Spuštění navigování k požadované euroklíč lokaci.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Spuštění navigování k požadované lokaci.", 2) | Out-Null

$d.Content.Find.Execute("Zobrazení všech euroklíč lokací na mapě.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Zobrazení všech lokací na mapě.", 2) | Out-Null

# --- 5. " + Scénáře" -> " +" + " jejich" + " " + "s" + "cénáře" -----------
$find2 = $d.Content
$find2.Find.Execute(" + Scénáře", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
Replace-RangeWithRuns $d $find2.Start $find2.End @(" +", " jejich", " ", "s", "cénáře")
